$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 141, shifting rows 141:169 down to 142:170
$ws.Rows.Item(141).Insert()

# Populate the new row 141 with the new record's data
$ws.Cells.Item(141, 1).Value = 5
$ws.Cells.Item(141, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(141, 3).Value = "Maule"
$ws.Cells.Item(141, 4).Value = 44694
$ws.Cells.Item(141, 4).NumberFormat = $ws.Cells.Item(142, 4).NumberFormat
$ws.Cells.Item(141, 5).Value = 7
$ws.Cells.Item(141, 6).Value = 100112017
$ws.Cells.Item(141, 7).Value = "Apio"
$ws.Cells.Item(141, 8).Value = "Americana (o)"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 600
$ws.Cells.Item(141, 11).Value = 7000
$ws.Cells.Item(141, 12).Value = 7000
$ws.Cells.Item(141, 13).Value = 7000
$ws.Cells.Item(141, 14).Value = "`$/docena de matas"
$ws.Cells.Item(141, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(141, 16).Value = 1167
$ws.Cells.Item(141, 17).Value = 6
$ws.Cells.Item(141, 18).Value = "Hortaliza"
